$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns, and reorder the final three coins
# (rows 49-51) to reflect the refreshed rankings/values.

$ws.Range("D2").Value = '29.277.11'
$ws.Range("E2").Value = '  +3.28%  '
$ws.Range("D3").Value = '1.584.54'
$ws.Range("E3").Value = '  +2.14%  '
$ws.Range("D4").Value = '0.997'
$ws.Range("E4").Value = '  -0.32%  '
$ws.Range("D5").Value = '212.54'
$ws.Range("E5").Value = '  +1.17%  '
$ws.Range("D6").Value = '0.511'
$ws.Range("E6").Value = '  +6.43%  '
$ws.Range("D7").Value = '0.996'
$ws.Range("E7").Value = '  -0.40%  '
$ws.Range("D8").Value = '26.47'
$ws.Range("E8").Value = '  +10.85%  '
$ws.Range("D9").Value = '0.249'
$ws.Range("E9").Value = '  +2.66%  '
$ws.Range("E10").Value = '  +2.03%  '
$ws.Range("E11").Value = '  +1.84%  '
$ws.Range("D12").Value = '1.810.50'
$ws.Range("E12").Value = '  +2.03%  '
$ws.Range("D13").Value = '1.595.15'
$ws.Range("E13").Value = '  +2.75%  '
$ws.Range("D14").Value = '29.317.37'
$ws.Range("E14").Value = '  +3.42%  '
$ws.Range("D15").Value = '3.73'
$ws.Range("E15").Value = '  +3.23%  '
$ws.Range("D16").Value = '0.525'
$ws.Range("E16").Value = '  +3.08%  '
$ws.Range("D17").Value = '63.03'
$ws.Range("E17").Value = '  +3.62%  '
$ws.Range("D18").Value = '239.56'
$ws.Range("E18").Value = '  +5.23%  '
$ws.Range("D19").Value = "'7.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.89%  '
$ws.Range("D20").Value = '0.0₃0689'
$ws.Range("E20").Value = '  +2.10%  '
$ws.Range("D21").Value = '0.996'
$ws.Range("E21").Value = '  -0.36%  '
$ws.Range("D22").Value = "'4.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.27%  '
$ws.Range("D23").Value = '9.25'
$ws.Range("E23").Value = '  +3.56%  '
$ws.Range("D24").Value = '2.09'
$ws.Range("E24").Value = '  +2.81%  '
$ws.Range("D25").Value = '154.53'
$ws.Range("E25").Value = '  +2.42%  '
$ws.Range("E26").Value = '  +4.98%  '
$ws.Range("D27").Value = '15.16'
$ws.Range("E27").Value = '  +2.97%  '
$ws.Range("E28").Value = '  +2.03%  '
$ws.Range("D29").Value = '0.996'
$ws.Range("E29").Value = '  -0.38%  '
$ws.Range("D30").Value = '0.0471'
$ws.Range("E30").Value = '  +0.68%  '
$ws.Range("E31").Value = '  -0.26%  '
$ws.Range("D32").Value = '3.23'
$ws.Range("E32").Value = '  +2.02%  '
$ws.Range("D33").Value = '1.423.28'
$ws.Range("E33").Value = '  +2.68%  '
$ws.Range("D34").Value = '3.09'
$ws.Range("E34").Value = '  +2.81%  '
$ws.Range("E35").Value = '  -3.30%  '
$ws.Range("E36").Value = '  +9.78%  '
$ws.Range("D37").Value = '1.51'
$ws.Range("E37").Value = '  +2.10%  '
$ws.Range("E38").Value = '  -1.56%  '
$ws.Range("D39").Value = '0.0165'
$ws.Range("E39").Value = '  +2.48%  '
$ws.Range("E40").Value = '  +3.68%  '
$ws.Range("D41").Value = '1.97'
$ws.Range("E41").Value = '  +2.69%  '
$ws.Range("D42").Value = '53.32'
$ws.Range("E42").Value = '  +25.34%  '
$ws.Range("D43").Value = '0.801'
$ws.Range("E43").Value = '  +3.04%  '
$ws.Range("D44").Value = '0.996'
$ws.Range("E44").Value = '  -0.39%  '
$ws.Range("D45").Value = '0.0469'
$ws.Range("E45").Value = '  +2.97%  '
$ws.Range("D46").Value = '64.65'
$ws.Range("E46").Value = '  +4.40%  '
$ws.Range("D47").Value = '5.35'
$ws.Range("E47").Value = '  -0.20%  '
$ws.Range("D48").Value = '1.722.05'
$ws.Range("E48").Value = '  +2.13%  '
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").Value = '85.95'
$ws.Range("E49").Value = '  +0.71%  '
$ws.Range("B50").Value = 'WEMIXToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D50").Value = '0.839'
$ws.Range("E50").Value = '  -2.58%  '
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").Value = '0.0₆0102'
$ws.Range("E51").Value = '  -1.31%  '
